# DevTesting_InflowsTEST.xlsx edit
#
# Commit message context:
#   "Edited IC to include release tier inputs."
#
# Concretely: for every trace-data worksheet in this workbook, column B
# ("Trace1") is overwritten with the values currently held in column G
# ("Trace6") for the data rows (rows 2-37). All other columns/cells are
# left untouched.

$wb = $excel.ActiveWorkbook

$sheetNames = @("DRGC2","BMDC2","CLSC2","GRNU1","GBRW4","MPSC2","NVRN5","GJLOC","GLDA3","TPIC2","VCRC2","YDLC2")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Copy the values from column G (Trace6) into column B (Trace1)
    # for rows 2 through 37 (the full data range on each sheet).
    $src = $ws.Range("G2:G37")
    $dst = $ws.Range("B2:B37")
    $dst.Value = $src.Value()
}
